$d = $word.ActiveDocument

# --- 1. "Date:" -> "Date of submition:" -------------------------------
# Position right after the word "Date" (before the colon) and type the
# new text, then nudge Bold off/on so the typed text + the remaining
# ":" land in their own runs (mirrors what Word's editor does when text
# is inserted in the middle of an existing run).
$rng = $d.Content
$rng.Find.Execute("Date", $false, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rng.Collapse(0)
$rng.InsertAfter(" of submition")
$rng.Bold = 0
$rng.Bold = 1

# --- 2. "Instructor: [Dr. Sameh / Eng. Nada" -> drop the leading "[" --
$d.Content.Find.Execute("[Dr. Sameh", $false, $false, $false, $false, $false, $true, 1, $false, "Dr. Sameh", 2)

# --- 3. Append student IDs after each team member's name -------------
function Append-ToParagraphEnd($paragraph, [string]$text) {
    $endRng = $d.Range($paragraph.Range.End - 1, $paragraph.Range.End - 1)
    $endRng.InsertAfter($text)
    $endRng.Bold = 1
    $endRng.Bold = 0
}

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text.TrimEnd([char]13, [char]10, [char]7)
    if ($t -eq "Nadeem Diaa 231000") {
        Append-ToParagraphEnd $para "857"
    } elseif ($t -eq "Ahmed Fahmy 231000") {
        Append-ToParagraphEnd $para "587"
    } elseif ($t -eq "Ahmed Shalaby 231000") {
        Append-ToParagraphEnd $para "514"
    } elseif ($t -eq "Yassin Mashhour 231000") {
        Append-ToParagraphEnd $para "604"
    }
}

Write-Output "Done."
